$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("summary_statistics")
$ws.Range("B2").Value = 347
$ws.Range("C2").Value = -9.210000000000001
$ws.Range("D2").Value = 15.55
$ws.Range("E2").Value = 1.3
$ws.Range("F2").Value = 0.6899999999999999
$ws.Range("G2").Value = 2.4
$ws.Range("J2").Value = 2.340215791392652
$ws.Range("B3").Value = 347
$ws.Range("B4").Value = 347
$ws.Range("B5").Value = 347
$ws.Range("B6").Value = 347
$ws.Range("E6").Value = 0.14
$ws.Range("G6").Value = 0.35
$ws.Range("B7").Value = 347
$ws.Range("E7").Value = 0.4
$ws.Range("B8").Value = 347
$ws.Range("E8").Value = 0.3
$ws.Range("G8").Value = 0.46
$ws.Range("B9").Value = 347
$ws.Range("E9").Value = 0.31
$ws.Range("G9").Value = 0.46
$ws.Range("B10").Value = 347
$ws.Range("E10").Value = 0.25
$ws.Range("H10").Value = 0.5
$ws.Range("J10").Value = 0.5
$ws.Range("B11").Value = 347
$ws.Range("E11").Value = 0.04
$ws.Range("G11").Value = 0.2
$ws.Range("B12").Value = 347
$ws.Range("E12").Value = 0.32
$ws.Range("G12").Value = 0.47
$ws.Range("B13").Value = 347
$ws.Range("E13").Value = 824.75
$ws.Range("F13").Value = 188.52
$ws.Range("G13").Value = 2360.51
$ws.Range("H13").Value = 503.29
$ws.Range("I13").Value = 26.04766666666667
$ws.Range("B14").Value = 343
$ws.Range("D14").Value = 97.90000000000001
$ws.Range("G14").Value = 8.859999999999999
$ws.Range("B15").Value = 343
$ws.Range("D15").Value = 55.1
$ws.Range("E15").Value = 4.03
$ws.Range("G15").Value = 7.75
$ws.Range("H15").Value = 3.2
$ws.Range("J15").Value = 3.6
$ws.Range("B16").Value = 343
$ws.Range("D16").Value = 97.90000000000001
$ws.Range("E16").Value = 13.24
$ws.Range("F16").Value = 3.5
$ws.Range("G16").Value = 21.49
$ws.Range("H16").Value = 14.9
$ws.Range("J16").Value = 15.6
$ws.Range("B17").Value = 343
$ws.Range("E17").Value = 25.57
$ws.Range("F17").Value = 14.8
$ws.Range("G17").Value = 26.72
$ws.Range("H17").Value = 30.7
$ws.Range("I17").Value = 6.5
$ws.Range("J17").Value = 37.2
$ws.Range("B18").Value = 343
$ws.Range("D18").Value = 12.4
$ws.Range("E18").Value = 0.26
$ws.Range("G18").Value = 0.82
$ws.Range("H18").Value = 0.3
$ws.Range("J18").Value = 0.3
$ws.Range("B19").Value = 343
$ws.Range("D19").Value = 2.7
$ws.Range("G19").Value = 0.16
$ws.Range("B20").Value = 343
$ws.Range("D20").Value = 33.6
$ws.Range("E20").Value = 4.11
$ws.Range("F20").Value = 3.5
$ws.Range("G20").Value = 3.44
$ws.Range("H20").Value = 4.2
$ws.Range("I20").Value = 1.5
$ws.Range("J20").Value = 5.7
$ws.Range("B21").Value = 343
$ws.Range("D21").Value = 98.90000000000001
$ws.Range("E21").Value = 51.21
$ws.Range("F21").Value = 55.8
$ws.Range("G21").Value = 31.75
$ws.Range("H21").Value = 57.3
$ws.Range("I21").Value = 23.2
$ws.Range("J21").Value = 80.5
$ws.Range("A22").Value = "percentfreelunchqualified"
$ws.Range("B22").Value = 306
$ws.Range("C22").Value = 0.01
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 0.47
$ws.Range("F22").Value = 0.42
$ws.Range("G22").Value = 0.28
$ws.Range("H22").Value = 0.46
$ws.Range("I22").Value = 0.2275368694221153
$ws.Range("J22").Value = 0.6908896846164441
$ws.Range("B23").Value = 347
$ws.Range("E23").Value = 51.39
$ws.Range("F23").Value = 53.22
$ws.Range("G23").Value = 27.91
$ws.Range("H23").Value = 48.42
$ws.Range("I23").Value = 26.96
$ws.Range("J23").Value = 75.38

$ws = $wb.Worksheets.Item("dichotomous_stats")
$ws.Range("B2").Value = 307
$ws.Range("C2").Value = 40
$ws.Range("D2").Value = 0.401
$ws.Range("E2").Value = 1.345
$ws.Range("F2").Value = 0.944
$ws.Range("G2").Value = 1.435
$ws.Range("H2").Value = 0.156
$ws.Range("I2").Value = 69.881
$ws.Range("J2").Value = -0.156
$ws.Range("K2").Value = 0.958
$ws.Range("B3").Value = 258
$ws.Range("C3").Value = 89
$ws.Range("D3").Value = 0.413
$ws.Range("E3").Value = 1.405
$ws.Range("F3").Value = 0.991
$ws.Range("G3").Value = 1.55
$ws.Range("H3").Value = 0.123
$ws.Range("I3").Value = 186.08
$ws.Range("J3").Value = -0.113
$ws.Range("K3").Value = 0.9399999999999999
$ws.Range("B4").Value = 257
$ws.Range("C4").Value = 90
$ws.Range("D4").Value = 0.5659999999999999
$ws.Range("E4").Value = 1.445
$ws.Range("F4").Value = 0.88
$ws.Range("G4").Value = 2.308
$ws.Range("H4").Value = 0.022
$ws.Range("I4").Value = 228.052
$ws.Range("J4").Value = 0.083
$ws.Range("K4").Value = 1.049
$ws.Range("B5").Value = 297
$ws.Range("C5").Value = 50
$ws.Range("D5").Value = 0.435
$ws.Range("E5").Value = 1.361
$ws.Range("F5").Value = 0.927
$ws.Range("G5").Value = 1.336
$ws.Range("H5").Value = 0.186
$ws.Range("I5").Value = 74.14700000000001
$ws.Range("J5").Value = -0.214
$ws.Range("K5").Value = 1.083
$ws.Range("B6").Value = 207
$ws.Range("C6").Value = 140
$ws.Range("D6").Value = 0.367
$ws.Range("E6").Value = 1.447
$ws.Range("F6").Value = 1.08
$ws.Range("G6").Value = 1.463
$ws.Range("H6").Value = 0.145
$ws.Range("I6").Value = 336.817
$ws.Range("J6").Value = -0.126
$ws.Range("K6").Value = 0.86
$ws.Range("B7").Value = 243
$ws.Range("C7").Value = 104
$ws.Range("D7").Value = 0.361
$ws.Range("E7").Value = 1.407
$ws.Range("F7").Value = 1.046
$ws.Range("G7").Value = 1.488
$ws.Range("H7").Value = 0.138
$ws.Range("I7").Value = 278.122
$ws.Range("J7").Value = -0.117
$ws.Range("K7").Value = 0.839
$ws.Range("B8").Value = 238
$ws.Range("C8").Value = 109
$ws.Range("D8").Value = 0.466
$ws.Range("F8").Value = 0.979
$ws.Range("G8").Value = 1.893
$ws.Range("H8").Value = 0.059
$ws.Range("I8").Value = 280.658
$ws.Range("J8").Value = -0.019
$ws.Range("K8").Value = 0.952
$ws.Range("B9").Value = 260
$ws.Range("C9").Value = 87
$ws.Range("D9").Value = 0.157
$ws.Range("E9").Value = 1.338
$ws.Range("F9").Value = 1.181
$ws.Range("G9").Value = 0.629
$ws.Range("H9").Value = 0.53
$ws.Range("I9").Value = 212.202
$ws.Range("J9").Value = -0.335
$ws.Range("K9").Value = 0.649
$ws.Range("B10").Value = 332
$ws.Range("D10").Value = 0.487
$ws.Range("E10").Value = 1.32
$ws.Range("F10").Value = 0.833
$ws.Range("G10").Value = 1.064
$ws.Range("H10").Value = 0.302
$ws.Range("I10").Value = 16.701
$ws.Range("J10").Value = -0.48
$ws.Range("K10").Value = 1.454
$ws.Range("B11").Value = 235
$ws.Range("C11").Value = 112
$ws.Range("D11").Value = 0.645
$ws.Range("E11").Value = 1.507
$ws.Range("F11").Value = 0.862
$ws.Range("G11").Value = 2.731
$ws.Range("H11").Value = 0.007
$ws.Range("I11").Value = 314.803
$ws.Range("J11").Value = 0.18
$ws.Range("K11").Value = 1.11

$ws = $wb.Worksheets.Item("anovas")
$ws.Range("C2").Value = 128.0615123297667
$ws.Range("D2").Value = 42.68717077658889
$ws.Range("E2").Value = 7.855281939139749
$ws.Range("F2").Value = 0.00004392986621105927
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 19.19273953776488
$ws.Range("D3").Value = 6.397579845921626
$ws.Range("E3").Value = 1.112312803280862
$ws.Range("F3").Value = 0.3440896756659662
$ws.Range("C4").Value = 411.755067866272
$ws.Range("D4").Value = 21.67131936138274
$ws.Range("E4").Value = 4.484467639607457
$ws.Range("F4").Value = 0.000000004727438165589188

$ws = $wb.Worksheets.Item("continuous_correlations")
$ws.Range("B2").Value = 0.014
$ws.Range("C2").Value = 0.253
$ws.Range("D2").Value = 0.801
$ws.Range("E2").Value = 345
$ws.Range("F2").Value = -0.092
$ws.Range("G2").Value = 0.119
$ws.Range("B3").Value = -0.04
$ws.Range("C3").Value = -0.735
$ws.Range("D3").Value = 0.463
$ws.Range("E3").Value = 341
$ws.Range("F3").Value = -0.145
$ws.Range("G3").Value = 0.066
$ws.Range("B4").Value = 0.018
$ws.Range("C4").Value = 0.336
$ws.Range("D4").Value = 0.737
$ws.Range("E4").Value = 341
$ws.Range("F4").Value = -0.08799999999999999
$ws.Range("G4").Value = 0.124
$ws.Range("B5").Value = -0.007
$ws.Range("C5").Value = -0.138
$ws.Range("D5").Value = 0.89
$ws.Range("E5").Value = 341
$ws.Range("F5").Value = -0.113
$ws.Range("G5").Value = 0.098
$ws.Range("B6").Value = 0.1
$ws.Range("C6").Value = 1.85
$ws.Range("D6").Value = 0.065
$ws.Range("E6").Value = 341
$ws.Range("F6").Value = -0.006
$ws.Range("G6").Value = 0.203
$ws.Range("B7").Value = -0.035
$ws.Range("C7").Value = -0.649
$ws.Range("D7").Value = 0.517
$ws.Range("E7").Value = 341
$ws.Range("F7").Value = -0.141
$ws.Range("G7").Value = 0.07099999999999999
$ws.Range("B8").Value = -0.06900000000000001
$ws.Range("C8").Value = -1.284
$ws.Range("D8").Value = 0.2
$ws.Range("E8").Value = 341
$ws.Range("F8").Value = -0.174
$ws.Range("G8").Value = 0.037
$ws.Range("B9").Value = -0.019
$ws.Range("C9").Value = -0.343
$ws.Range("D9").Value = 0.732
$ws.Range("E9").Value = 341
$ws.Range("F9").Value = -0.124
$ws.Range("G9").Value = 0.08699999999999999
$ws.Range("B10").Value = -0.06900000000000001
$ws.Range("C10").Value = -1.279
$ws.Range("D10").Value = 0.202
$ws.Range("E10").Value = 341
$ws.Range("F10").Value = -0.174
$ws.Range("G10").Value = 0.037
$ws.Range("A11").Value = "percentfreelunchqualified"
$ws.Range("B11").Value = 0.001
$ws.Range("C11").Value = 0.017
$ws.Range("D11").Value = 0.987
$ws.Range("E11").Value = 304
$ws.Range("F11").Value = -0.111
$ws.Range("G11").Value = 0.113
$ws.Range("B12").Value = 0.03
$ws.Range("C12").Value = 0.552
$ws.Range("D12").Value = 0.581
$ws.Range("E12").Value = 345
$ws.Range("F12").Value = -0.076
$ws.Range("G12").Value = 0.135

Write-Output "done"